# Applies the "partial updates based on comments from JN" edit:
#  - Notes Master date placeholder cache: 2/19/2025 -> 2/20/2025
#  - "First Layer: Hash" -> "First Step: Hash"
#  - "Second Layer: Hash & Salt" -> "Second Step: Hash & Salt"
#  - "Third Layer: Encryption" -> "Third Step: " + "Encryption" (two runs)
#  - "API Request Encryption" -> "API Request Application-Level Encryption"

$p = $ppt.ActivePresentation

# --- Notes Master date field (auto-date cache bump by one day) ---
$nm = $p.NotesMaster
$hf = $nm.HeadersFooters
$hf.DateAndTime.Text = "2/20/2025"

# --- Slide 1 ---
$s = $p.Slides.Item(1)

# "First Layer: Hash" -> "First Step: Hash"
$shFirst = $s.Shapes.Item("Rectangle 5")
$shFirst.TextFrame.TextRange.Paragraphs(1, 1).Runs(1, 1).Text = "First Step: Hash"

# "Second Layer: Hash & Salt" -> "Second Step: Hash & Salt"
$shSecond = $s.Shapes.Item("Rectangle 7")
$shSecond.TextFrame.TextRange.Paragraphs(1, 1).Runs(1, 1).Text = "Second Step: Hash & Salt"

# "Third Layer: Encryption" -> two runs: "Third Step: " + "Encryption"
$shThird = $s.Shapes.Item("Rectangle 9")
$thirdRun = $shThird.TextFrame.TextRange.Paragraphs(1, 1).Runs(1, 1)
$thirdRun.InsertBefore("Third Step: ") | Out-Null
$shThird.TextFrame.TextRange.Paragraphs(1, 1).Runs(2, 1).Text = "Encryption"

# "API Request Encryption" -> "API Request Application-Level Encryption"
$shApi = $s.Shapes.Item("Rectangle 13")
$shApi.TextFrame.TextRange.Paragraphs(1, 1).Runs(2, 1).Text = "API Request Application-Level Encryption"
